$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 9789.963
$ws.Range("I70").Value = 5237.5
$ws.Range("K70").Value = 15712.5
$ws.Range("M70").Value = -15442.5
$ws.Range("H73").Value = 9789.963
$ws.Range("I73").Value = 5237.5
$ws.Range("K73").Value = 15712.5
$ws.Range("M73").Value = -14776.5
$ws.Range("H92").Value = 1354.12
$ws.Range("I92").Value = 1056.3684
$ws.Range("J92").Value = 2297
$ws.Range("K92").Value = 1056.3684
$ws.Range("L92").Value = 2297
$ws.Range("M92").Value = 191.6315999999999
$ws.Range("N92").Value = -4793
$ws.Range("H107").Value = 1591.1282
$ws.Range("I107").Value = 1097.7727
$ws.Range("K107").Value = 1097.7727
$ws.Range("M107").Value = 822.2273
$ws.Range("H137").Value = 2604.303
$ws.Range("I137").Value = 2672.7727
$ws.Range("J137").Value = 2467.3635
$ws.Range("K137").Value = 8018.3181
$ws.Range("L137").Value = 7402.0905
$ws.Range("M137").Value = -5468.3181
$ws.Range("N137").Value = -12502.0905
$ws.Range("H138").Value = 2243.61
$ws.Range("J138").Value = 2486.3735
$ws.Range("L138").Value = 7459.120500000001
$ws.Range("N138").Value = -17739.1205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4612.9
$ws.Range("I45").Value = 3447.1428
$ws.Range("K45").Value = 3447.1428
$ws.Range("M45").Value = -3070.1428
$ws.Range("H80").Value = 65537.86
$ws.Range("J80").Value = 65537.86
$ws.Range("L80").Value = 65537.86
$ws.Range("N80").Value = -67533.86
$ws.Range("H83").Value = 65537.86
$ws.Range("J83").Value = 65537.86
$ws.Range("L83").Value = 196613.58
$ws.Range("N83").Value = -206597.58
$ws.Range("H97").Value = 4138.647
$ws.Range("I97").Value = 1409.4286
$ws.Range("J97").Value = 16875
$ws.Range("K97").Value = 1409.4286
$ws.Range("L97").Value = 16875
$ws.Range("M97").Value = -913.4286
$ws.Range("N97").Value = -17867
$ws.Range("H122").Value = 4650.5864
$ws.Range("J122").Value = 5482.1665
$ws.Range("L122").Value = 16446.4995
$ws.Range("N122").Value = -21346.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 450.83334
$ws.Range("I22").Value = 450.83334
$ws.Range("K22").Value = 450.83334
$ws.Range("M22").Value = -277.83334
$ws.Range("H86").Value = 2303011.8
$ws.Range("I86").Value = 3706381.5
$ws.Range("J86").Value = 6588.727
$ws.Range("K86").Value = 3706381.5
$ws.Range("L86").Value = 6588.727
$ws.Range("M86").Value = -3705258.5
$ws.Range("N86").Value = -8834.726999999999
$ws.Range("H89").Value = 2303011.8
$ws.Range("I89").Value = 3706381.5
$ws.Range("J89").Value = 6588.727
$ws.Range("K89").Value = 18531907.5
$ws.Range("L89").Value = 32943.635
$ws.Range("M89").Value = -18526291.5
$ws.Range("N89").Value = -44175.635
$ws.Range("H99").Value = 3889.6
$ws.Range("I99").Value = 3312.125
$ws.Range("K99").Value = 3312.125
$ws.Range("M99").Value = -1814.125
$ws.Range("H134").Value = 8020.255
$ws.Range("I134").Value = 4648.1875
$ws.Range("K134").Value = 13944.5625
$ws.Range("M134").Value = -11409.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3417.5483
$ws.Range("I31").Value = 2460.25
$ws.Range("J31").Value = 4022.158
$ws.Range("K31").Value = 2460.25
$ws.Range("L31").Value = 4022.158
$ws.Range("M31").Value = -2165.25
$ws.Range("N31").Value = -4612.157999999999
$ws.Range("H34").Value = 3417.5483
$ws.Range("I34").Value = 2460.25
$ws.Range("J34").Value = 4022.158
$ws.Range("K34").Value = 2460.25
$ws.Range("L34").Value = 4022.158
$ws.Range("M34").Value = -2258.25
$ws.Range("N34").Value = -4426.157999999999
$ws.Range("H53").Value = 50153.4
$ws.Range("J53").Value = 50153.4
$ws.Range("L53").Value = 50153.4
$ws.Range("N53").Value = -51367.4
$ws.Range("H58").Value = 2404.5625
$ws.Range("I58").Value = 1579.5
$ws.Range("K58").Value = 1579.5
$ws.Range("M58").Value = -1376.5
$ws.Range("H134").Value = 4549.075
$ws.Range("I134").Value = 1073.15
$ws.Range("J134").Value = 8025
$ws.Range("K134").Value = 3219.45
$ws.Range("L134").Value = 24075
$ws.Range("M134").Value = -684.4500000000003
$ws.Range("N134").Value = -29145
$ws.Range("H136").Value = 2404.5625
$ws.Range("I136").Value = 1579.5
$ws.Range("K136").Value = 4738.5
$ws.Range("M136").Value = -2188.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 59395.11
$ws.Range("J75").Value = 87875.086
$ws.Range("L75").Value = 263625.258
$ws.Range("N75").Value = -265621.258
$ws.Range("H78").Value = 59395.11
$ws.Range("J78").Value = 87875.086
$ws.Range("L78").Value = 790875.774
$ws.Range("N78").Value = -800859.774
$ws.Range("H129").Value = 51286.35
$ws.Range("I129").Value = 91546.82000000001
$ws.Range("K129").Value = 274640.46
$ws.Range("M129").Value = -269640.46
$ws.Range("H131").Value = 556769.25
$ws.Range("I131").Value = 1250544.5
$ws.Range("K131").Value = 3751633.5
$ws.Range("M131").Value = -3746593.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 44444
$ws.Range("J52").Value = 44444
$ws.Range("L52").Value = 44444
$ws.Range("N52").Value = -44962
$ws.Range("H80").Value = 3431.2307
$ws.Range("J80").Value = 3651.5625
$ws.Range("L80").Value = 3651.5625
$ws.Range("N80").Value = -5647.5625
$ws.Range("H83").Value = 3431.2307
$ws.Range("J83").Value = 3651.5625
$ws.Range("L83").Value = 18257.8125
$ws.Range("N83").Value = -28241.8125
$ws.Range("H102").Value = 20786.346
$ws.Range("I102").Value = 3463.7368
$ws.Range("K102").Value = 3463.7368
$ws.Range("M102").Value = -1841.7368
$ws.Range("H122").Value = 1410.7858
$ws.Range("I122").Value = 1240.3
$ws.Range("J122").Value = 1837
$ws.Range("K122").Value = 3720.9
$ws.Range("L122").Value = 5511
$ws.Range("M122").Value = -1270.9
$ws.Range("N122").Value = -10411

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2545.1667
$ws.Range("I22").Value = 493.5
$ws.Range("J22").Value = 3131.3572
$ws.Range("K22").Value = 493.5
$ws.Range("L22").Value = 3131.3572
$ws.Range("M22").Value = -198.5
$ws.Range("N22").Value = -3721.3572
$ws.Range("H27").Value = 2545.1667
$ws.Range("I27").Value = 493.5
$ws.Range("J27").Value = 3131.3572
$ws.Range("K27").Value = 493.5
$ws.Range("L27").Value = 3131.3572
$ws.Range("M27").Value = -386.5
$ws.Range("N27").Value = -3345.3572
$ws.Range("H40").Value = 5832.913
$ws.Range("I40").Value = 5676.0557
$ws.Range("K40").Value = 5676.0557
$ws.Range("M40").Value = -5540.0557
$ws.Range("H68").Value = 3486
$ws.Range("I68").Value = 3453
$ws.Range("K68").Value = 3453
$ws.Range("M68").Value = -2704
$ws.Range("H71").Value = 3486
$ws.Range("I71").Value = 3453
$ws.Range("K71").Value = 17265
$ws.Range("M71").Value = -13521
$ws.Range("H82").Value = 5996.913
$ws.Range("I82").Value = 9365.083000000001
$ws.Range("K82").Value = 9365.083000000001
$ws.Range("M82").Value = -9004.083000000001
$ws.Range("H85").Value = 5996.913
$ws.Range("I85").Value = 9365.083000000001
$ws.Range("K85").Value = 9365.083000000001
$ws.Range("M85").Value = -8117.083000000001
$ws.Range("H93").Value = 3411.5
$ws.Range("I93").Value = 1823
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 1823
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -575
$ws.Range("N93").Value = -7496
$ws.Range("H122").Value = 7141.923
$ws.Range("I122").Value = 6519.375
$ws.Range("K122").Value = 19558.125
$ws.Range("M122").Value = -17108.125
$ws.Range("H136").Value = 5564.36
$ws.Range("I136").Value = 4666.8887
$ws.Range("K136").Value = 14000.6661
$ws.Range("M136").Value = -11450.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2407.8928
$ws.Range("I122").Value = 2338.45
$ws.Range("K122").Value = 7015.349999999999
$ws.Range("M122").Value = -4565.349999999999
$ws.Range("H132").Value = 1399.7
$ws.Range("I132").Value = 1221.0714
$ws.Range("J132").Value = 1816.5
$ws.Range("K132").Value = 3663.2142
$ws.Range("L132").Value = 5449.5
$ws.Range("M132").Value = -1133.2142
$ws.Range("N132").Value = -10509.5
$ws.Range("H136").Value = 9037.173000000001
$ws.Range("I136").Value = 13355.883
$ws.Range("K136").Value = 40067.649
$ws.Range("M136").Value = -37517.649
